$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AD1:AF1, copying the
# formatting (bold, border, centered) used by the other header cells.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (W/L/T) for every player row.
$ws.Range("AD2:AD49").Value = 87
$ws.Range("AE2:AE49").Value = 75
$ws.Range("AF2:AF49").Value = 0
